# C5-PowerPoint.pptx edit
#  1) Slide 6's table switches to a different built-in table style.
#  2) The presentation's applied theme colour scheme switches from the
#     "Integral" palette to the stock "Office Theme" palette (this is the
#     effect of picking a different theme from the Design tab).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 -------------------------------------------
$slide  = $p.Slides.Item(6)
$shape  = $slide.Shapes.Item(2)
$table  = $shape.Table
$table.ApplyStyle("{4E3E0483-133E-435F-89E4-8FC8D0FEED87}")

# --- 2) Swap the applied theme's colour scheme ----------------------------
$theme  = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# index -> (scheme slot, target RGB as 0xBBGGRR COM colour value)
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
